# Phase-2 Add Frontend Agent Flow-Complete
#
# Updates the ticket-number column (Y) on the "NFTRTickets" sheet with a
# fresh batch of ticket numbers. The values look numeric but must stay
# text (they carry significant leading zeros), so the target cells are
# explicitly formatted as Text ("@") before the value is written - this
# mirrors the "Text" number format (numFmtId 49) that this workbook
# already uses elsewhere (e.g. LoginCredentials!F1) for the same kind of
# numeric-looking identifier.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFTRTickets")

$updates = @{
    "Y2" = "040820001136"
    "Y3" = "040820001132"
    "Y4" = "040820001137"
    "Y5" = "040820001138"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
